$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update existing figures with restated (smaller-unit) values
$ws.Range("D2").Value = 4077
$ws.Range("E2").Value = 153
$ws.Range("F2").Value = 153
$ws.Range("G2").Value = 117
$ws.Range("H2").Value = 73
$ws.Range("I2").Value = 35
$ws.Range("J2").Value = 38
$ws.Range("K2").Value = 4313
$ws.Range("L2").Value = 1612
$ws.Range("M2").Value = 2701
$ws.Range("N2").Value = 1943
$ws.Range("O2").Value = 758
$ws.Range("P2").Value = 146
$ws.Range("Q2").Value = 241
$ws.Range("R2").Value = -210
$ws.Range("S2").Value = -94
$ws.Range("T2").Value = 478
$ws.Range("U2").Value = -237
$ws.Range("V2").Value = 818
$ws.Range("W2").Value = 3.74
$ws.Range("X2").Value = 1.78
$ws.Range("Y2").Value = 1.74
$ws.Range("Z2").Value = 1.67
$ws.Range("AA2").Value = 59.69
$ws.Range("AB2").Value = 1281.03
$ws.Range("AC2").Value = 119
$ws.Range("AD2").Value = 48.71
$ws.Range("AE2").Value = 6647
$ws.Range("AF2").Value = 0.87
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 2.59
$ws.Range("AI2").Value = 125.99
$ws.Range("AJ2").Value = 29228750

# Row 3: update existing figures with restated (smaller-unit) values
$ws.Range("D3").Value = 3539
$ws.Range("E3").Value = -33
$ws.Range("F3").Value = -33
$ws.Range("G3").Value = -58
$ws.Range("H3").Value = -83
$ws.Range("I3").Value = -104
$ws.Range("J3").Value = 21
$ws.Range("K3").Value = 4670
$ws.Range("L3").Value = 2097
$ws.Range("M3").Value = 2574
$ws.Range("N3").Value = 1819
$ws.Range("O3").Value = 755
$ws.Range("P3").Value = 146
$ws.Range("Q3").Value = 231
$ws.Range("R3").Value = -495
$ws.Range("S3").Value = 132
$ws.Range("T3").Value = 435
$ws.Range("U3").Value = -204
$ws.Range("V3").Value = 1407
$ws.Range("W3").Value = -0.92
$ws.Range("X3").Value = -2.34
$ws.Range("Y3").Value = -5.54
$ws.Range("Z3").Value = -1.85
$ws.Range("AA3").Value = 81.47
$ws.Range("AB3").Value = 1212.83
$ws.Range("AC3").Value = -357
$ws.Range("AD3").Value = -16.21
$ws.Range("AE3").Value = 6223
$ws.Range("AF3").Value = 0.93
$ws.Range("AG3").Value = 180
$ws.Range("AH3").Value = 3.11
$ws.Range("AI3").Value = -50.47
$ws.Range("AJ3").Value = 29228750

# Row 4: update existing figures with restated (smaller-unit) values
$ws.Range("D4").Value = 4447
$ws.Range("E4").Value = 122
$ws.Range("F4").Value = 122
$ws.Range("G4").Value = 87
$ws.Range("H4").Value = 44
$ws.Range("I4").Value = 25
$ws.Range("J4").Value = 19
$ws.Range("K4").Value = 4984
$ws.Range("L4").Value = 2358
$ws.Range("M4").Value = 2626
$ws.Range("N4").Value = 1768
$ws.Range("O4").Value = 858
$ws.Range("P4").Value = 146
$ws.Range("Q4").Value = 174
$ws.Range("R4").Value = -266
$ws.Range("S4").Value = 253
$ws.Range("T4").Value = 835
$ws.Range("U4").Value = -661
$ws.Range("V4").Value = 1618
$ws.Range("W4").Value = 2.75
$ws.Range("X4").Value = 0.98
$ws.Range("Y4").Value = 1.38
$ws.Range("Z4").Value = 0.9
$ws.Range("AA4").Value = 89.81
$ws.Range("AB4").Value = 1170.86
$ws.Range("AC4").Value = 84
$ws.Range("AD4").Value = 56.93
$ws.Range("AE4").Value = 6047
$ws.Range("AF4").Value = 0.8
$ws.Range("AG4").Value = 200
$ws.Range("AH4").Value = 4.16
$ws.Range("AI4").Value = 236.71
$ws.Range("AJ4").Value = 29228750

# Row 5: update existing figures with restated (smaller-unit) values
$ws.Range("D5").Value = 5637
$ws.Range("E5").Value = 196
$ws.Range("F5").Value = 196
$ws.Range("G5").Value = 112
$ws.Range("H5").Value = 61
$ws.Range("I5").Value = 83
$ws.Range("J5").Value = -32
$ws.Range("K5").Value = 4762
$ws.Range("L5").Value = 2171
$ws.Range("M5").Value = 2591
$ws.Range("N5").Value = 1861
$ws.Range("O5").Value = 804
$ws.Range("P5").Value = 146
$ws.Range("Q5").Value = 399
$ws.Range("R5").Value = -113
$ws.Range("S5").Value = -330
$ws.Range("T5").Value = 153
$ws.Range("U5").Value = 246
$ws.Range("V5").Value = 1344
$ws.Range("W5").Value = 3.48
$ws.Range("X5").Value = 1.08
$ws.Range("Y5").Value = 4.56
$ws.Range("Z5").Value = 1.25
$ws.Range("AA5").Value = 83.76000000000001
$ws.Range("AB5").Value = 1253.28
$ws.Range("AC5").Value = 283
$ws.Range("AD5").Value = 18.37
$ws.Range("AE5").Value = 6366
$ws.Range("AF5").Value = 0.82
$ws.Range("AG5").Value = 200
$ws.Range("AH5").Value = 3.85
$ws.Range("AI5").Value = 70.66
$ws.Range("AJ5").Value = 29228750

# Row 6: update existing figures with restated (smaller-unit) values
$ws.Range("D6").Value = 5635
$ws.Range("E6").Value = 234
$ws.Range("F6").Value = 234
$ws.Range("G6").Value = 201
$ws.Range("H6").Value = 144
$ws.Range("I6").Value = 120
$ws.Range("K6").Value = 5427
$ws.Range("L6").Value = 2324
$ws.Range("M6").Value = 3103
$ws.Range("N6").Value = 2271
$ws.Range("P6").Value = 146
$ws.Range("Q6").Value = 183
$ws.Range("R6").Value = -141
$ws.Range("S6").Value = -136
$ws.Range("T6").Value = 208
$ws.Range("U6").Value = -25
$ws.Range("V6").Value = 1318
$ws.Range("W6").Value = 4.15
$ws.Range("X6").Value = 2.56
$ws.Range("Y6").Value = 5.82
$ws.Range("Z6").Value = 2.84
$ws.Range("AA6").Value = 74.90000000000001
$ws.Range("AB6").Value = 1312.13
$ws.Range("AC6").Value = 412
$ws.Range("AD6").Value = 13.56
$ws.Range("AE6").Value = 7855
$ws.Range("AF6").Value = 0.71
$ws.Range("AG6").Value = 220
$ws.Range("AH6").Value = 3.94
$ws.Range("AI6").Value = 52.92
$ws.Range("AJ6").Value = 29228750

# Row 7: remove forecast figures (estimate columns no longer populated)
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8: remove forecast figures (estimate columns no longer populated)
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9: remove forecast figures (estimate columns no longer populated)
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
